$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newAverage = 60.146050580527778

$pValues = @{
    2  = 0.99938827753067017
    3  = 0.99944943189620972
    4  = 0.99951064586639404
    5  = 0.99957180023193359
    6  = 0.99963295459747314
    7  = 0.99969416856765747
    8  = 0.99975532293319702
    9  = 0.99981647729873657
    10 = 0.99987763166427612
    11 = 0.99993884563446045
    12 = 0.99994492530822754
    13 = 0.9999510645866394
    14 = 0.99995720386505127
    15 = 0.99996328353881836
    16 = 0.99996942281723022
    17 = 0.99997550249099731
    18 = 0.99998164176940918
    19 = 0.99998778104782104
    20 = 0.99999386072158813
    21 = 0.99999451637268066
    22 = 0.99999511241912842
    23 = 0.99999570846557617
    24 = 0.99999630451202393
    25 = 0.99999696016311646
    26 = 0.99999755620956421
    27 = 0.99999815225601196
    28 = 0.99999874830245972
    29 = 0.99999940395355225
    30 = 0.99999946355819702
    31 = 0.9999995231628418
    32 = 0.99999958276748657
    33 = 0.99999964237213135
    34 = 0.99999970197677612
    35 = 0.9999997615814209
    36 = 0.99999982118606567
    37 = 0.99999988079071045
    38 = 0.99999994039535522
}

for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 3).Value = $newAverage
    $ws.Cells.Item($row, 4).Value = $pValues[$row]
}
